$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14 (A14): insert underlined URL run and add trailing space after "dong 1" ---
$b64_14 = "YwBoANceIABuAOAAeQAgAGwA4AAgAHgAdQDRHm4AZwAgADIAIAB2AOAAIAAgADMAIABoAHQAdABwAHMAOgAvAC8AdABoAGkAcwBpAHMAaQBtAGEAZwBlAHUAcgBsAC4AYwBvAG0AIABkAPIAbgBnAAoAZADyAG4AZwAgADEAIAAKAAoAZADyAG4AZwAgADIACgAKAAoAZADyAG4AZwAgADMA"
$bytes_14 = [Convert]::FromBase64String($b64_14)
$text_14 = [System.Text.Encoding]::Unicode.GetString($bytes_14)

$rng14 = $ws.Range("A14")
$rng14.Value = $text_14

# Reapply formatting for the run " 3" (bold + underline) -> chars 23-24
$rng14.Characters(23, 2).Font.Bold = $true
$rng14.Characters(23, 2).Font.Underline = $true

# New run: " https://thisisimageurl.com " (underline only) -> chars 25-52
$rng14.Characters(25, 28).Font.Underline = $true

# --- Row 16 (A16): prepend "111 " and insert URL into the -100% line ---
$b64_16 = "MQAxADEAIAACDkkOLQ4hDjkOJQ4qDjQOGQ4EDkkOMg5ADh4ONA5IDiEOQA4VDjQOIQ4KAAoALQAxADAAMAAlACAAaAB0AHQAcABzADoALwAvAHQAaABpAHMAaQBzAGkAbQBhAGcAZQB1AHIAbAAuAGMAbwBtACAAIABDAG8AbABvAHIAIABWAG8AbAB1AG0AZQAgAFEAdQBhAG4AdAB1AG0AIABEAG8AdAAgACoONQ4qDicOIg4qDhQOIAAxADAAMAAlAAoAUQB1AGEAbgB0AHUAbQAgAEgARABSACAAQA4cDiIOIw4yDiIOJQ4wDkAOLQ41DiIOFA4XDjUOSA4LDkgOLQ4ZDi0OIg44DkgOQw4ZDiAOMg4eDgoAQQBpAHIAUwBsAGkAbQAgABUOMQ4nDkAOBA4jDjcOSA4tDgcOGg4yDgcOIAAqDicOIg4HDjIOIQ4gABsOIw4wDisOIg4xDhQOHg43DkkOGQ4XDjUOSA4KAFMAbQBhAHIAdAAgAEgAdQBiACAABA4xDhQOKg4jDiMOBA4tDhkOQA4XDhkOFQ5MDhoOMQ4ZDkAOFw40DgcOIw4nDiEORA4nDkkOQw4ZDhcONQ5IDkAOFA41DiIOJw4gAAoABA44DhMOKg4hDhoOMQ4VDjQOQA4JDh4OMg4wDgoAQQ4aDiMOGQ4UDkwOIAA6ACAAUwBBAE0AUwBVAE4ARwAKAAsONQ4jDjUOKg5MDiAAOgAgAFEANgA1AEMACgAEDicOMg4hDioOOQ4HDiAAKAALDiEOLgApACAAOgAgADYANAAuADQAMQAKAAQOJw4yDiEOAQ4nDkkOMg4HDiAAKAALDiEOLgApACAAOgAgADEAMQAxAC4AOAAzAAoABA4nDjIOIQ4lDjYOAQ4gACgACw4hDi4AKQAgADoAIAAyAC4ANQA3AAoAGQ5JDjMOKw4ZDjEOAQ4gACgAAQ4BDi4AKQAgADoAIAAxADAALgA3AAoAAg4ZDjIOFA4rDhkOSQ4yDggOLQ4gACgAGQ40DkkOJw4pACAAOgAgADUAMAAiAAoABA4nDjIOIQ4lDjAOQA4tDjUOIg4UDisOGQ5JDjIOCA4tDiAAOgAgADQASwAgACgAMwAsADgANAAwACAAeAAgADIALAAxADYAMAApACAACgAKDhkONA4UDisOGQ5JDjIOCA4tDiAAOgAgAFEATABFAEQACgBEAEkARwBJAFQAQQBMACAAVABWACAAQgBVAEkATABUACAASQBOACAAKABZAEUAUwAvAE4ATwApACAAOgAgAFkARQBTAAoAUwBNAEEAUgBUACAAVABWACAAOgAgAFQASQBaAEUATgAKAEgARABSACAARgBPAFIATQBBAFQAIAA6ACAAUQB1AGEAbgB0AHUAbQAgAEgARABSAAoASABEAE0ASQAgACgACg5IDi0OBw4pACAAOgAgADMACgBVAFMAQgAgACgACg5IDi0OBw4pACAAOgAgADIACgAbDiMOMA5ADiAOFw4jDjUOQg4hDhcOIAA6ACAATwBOAEUAIABSAEUATQBPAFQARQAKAAQOJw4yDiEOJQ42DgEOIABUAFYAIAAjDicOIQ4CDjIOFQ4xDkkOBw4gACgACw4hDi4AKQAgADoAIAAxADkALgA5ADEACgAEDicOMg4hDioOOQ4HDiAAVABWACAAIw4nDiEOAg4yDhUOMQ5JDgcOIAAoAAsOIQ4uACkAIAA6ACAANwAwAC4AOQA0AAoAJw40DhgONQ5DDgoOSQ4HDjIOGQ4KAEMOCg5JDioOMw4rDiMOMQ4aDiMOMQ4aDgoOIQ5ADh4ONw5IDi0OBA4nDjIOIQ4aDjEOGQ5ADhcONA4HDgoABA4zDkEOGQ4wDhkOMw4KABYOLQ4UDhsOJQ4xDkoOAQ4qDjIOIg5EDh8OLQ4tDgEOCA4yDgEOQA4VDkkOMg4jDjEOGg5EDh8OHw5JDjIOFw41DkgOHA4ZDjEOBw4gAEEOJQ4wDkAOCg5HDhQOHA4lDjQOFQ4gDjEOEw4RDkwOFA5JDicOIg4cDkkOMg5BDisOSQ4HDkEOJQ4wDhkOOA5IDiEOIABADh4ONw5IDi0OGw5JDi0OBw4BDjEOGQ4jDi0OIg4CDjkOFA4CDjUOFA4KACgONg4BDikOMg4EDjkOSA4hDjcOLQ4BDjIOIw4VDjQOFA4VDjEOSQ4HDiAAQQ4lDjAOJw40DhgONQ4BDjIOIw5DDgoOSQ4HDjIOGQ4tDiIOSA4yDgcOJQ4wDkAOLQ41DiIOFA4KAAIOSQ4tDiEOOQ4lDgEOMg4jDiMOMQ4aDhsOIw4wDgEOMQ4ZDgoAAQ4yDiMOIw4xDhoOGw4jDjAOAQ4xDhkOIAAoABsONQ4pACAAOgAgADEACgACDkkOLQ4hDjkOJQ4VDjQOFA4VDkgOLQ4oDjkOGQ4iDkwOGg4jDjQOAQ4yDiMOCgBDAEEATABMACAAQwBFAE4AVABFAFIAIAAxADIAOAAyACAAGg4jDjQOKQ4xDhcOIABEDhcOIg4LDjEOIQ4LDjgOBw4gAC0ONA5ADiUOBA5CDhcOIw4ZDjQOBA4qDkwOIAAIDjMOAQ4xDhQOIAAoDjkOGQ4iDkwOGg4jDjQOAQ4yDiMOQA4bDjQOFA5DDisOSQ4aDiMONA4BDjIOIw4VDjEOSQ4HDkEOFQ5IDicOMQ4ZDggOMQ4ZDhcOIw5MDhYONg4HDicOMQ4ZDkAOKg4yDiMOTA4gADkALgAwADAAIAAxADcALgAwADAAIAAZDi4A"
$bytes_16 = [Convert]::FromBase64String($b64_16)
$text_16 = [System.Text.Encoding]::Unicode.GetString($bytes_16)

$ws.Range("A16").Value = $text_16
